$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.816.36"
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").Value = "2.928.25"
$ws.Range("E3").Value = "  +0.59%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "570.65"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.78%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "145.77"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.31%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "2.925.43"
$ws.Range("E8").Value = "  +0.56%  "
$ws.Range("E9").Value = "  -0.44%  "
$ws.Range("E10").Value = "  +0.00%  "
$ws.Range("E11").Value = "  +1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.436"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.28%  "
$ws.Range("E13").Value = "  +1.68%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "33.01"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.78%  "
$ws.Range("E15").Value = "  +0.76%  "
$ws.Range("D16").Value = "3.413.33"
$ws.Range("E16").Value = "  +0.54%  "
$ws.Range("D17").Value = "62.667.83"
$ws.Range("E17").Value = "  +1.32%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "6.66"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +1.08%  "
$ws.Range("D19").Value = "2.922.07"
$ws.Range("E19").Value = "  +0.38%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "434.22"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.64%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.19"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.78%  "
$ws.Range("E22").Value = "  +0.78%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.92"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.13%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "79.00"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.29%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "12.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +1.04%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "10.09"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -1.14%  "
$ws.Range("E27").Value = "  +0.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.04"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.29%  "
$ws.Range("E29").Value = "  +6.03%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "7.07"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.74%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "2.51"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.91%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.04"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -2.88%  "
$ws.Range("E33").Value = "  -1.15%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.97"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.47%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.998"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.22%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.959"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.77%  "
$ws.Range("E37").Value = "  -1.15%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.97"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.61%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "49.09"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.01%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.92"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -4.00%  "
$ws.Range("E41").Value = "  -0.70%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "41.28"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.13"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.58%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.270"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.70%  "
$ws.Range("D45").Value = "2.716.57"
$ws.Range("E45").Value = "  +0.72%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0341"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +1.26%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "133.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "354.72"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +3.22%  "
$ws.Range("E50").Value = "  +15.86%  "
$ws.Range("E51").Value = "  -0.55%  "
